$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl11"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.842606
$ws.Range("H2").Value = 2.527818
$ws.Range("I2").Value = 0.0108780433452729
$ws.Range("J2").Value = 0.0108780433452729
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 1.863797
$ws.Range("N2").Value = 5.591391
$ws.Range("O2").Value = 0.5455000708290748
$ws.Range("P2").Value = 0.5455000708290748
$ws.Range("Q2").Value = 1.570446534982
$ws.Range("R2").Value = 14.134018814838
$ws.Range("S2").Value = 0.005933973415328114
$ws.Range("T2").Value = 0.005933973415328114

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl11"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.842606
$ws.Range("H3").Value = 2.527818
$ws.Range("I3").Value = 0.0108780433452729
$ws.Range("J3").Value = 0.0108780433452729
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.552879
$ws.Range("N3").Value = 4.658637000000001
$ws.Range("O3").Value = 0.4544999291709252
$ws.Range("P3").Value = 0.4544999291709252
$ws.Range("Q3").Value = 1.308465162674
$ws.Range("R3").Value = 11.776186464066
$ws.Range("S3").Value = 0.004944069929944789
$ws.Range("T3").Value = 0.004944069929944789

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl11"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 69.05064766666666
$ws.Range("H4").Value = 207.151943
$ws.Range("I4").Value = 0.8914438519749055
$ws.Range("J4").Value = 0.8914438519749054
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 1.863797
$ws.Range("N4").Value = 5.591391
$ws.Range("O4").Value = 0.5455000708290748
$ws.Range("P4").Value = 0.5455000708290748
$ws.Range("Q4").Value = 128.6963899691903
$ws.Range("R4").Value = 1158.267509722713
$ws.Range("S4").Value = 0.4862826843924543
$ws.Range("T4").Value = 0.4862826843924541

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl11"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 69.05064766666666
$ws.Range("H5").Value = 207.151943
$ws.Range("I5").Value = 0.8914438519749055
$ws.Range("J5").Value = 0.8914438519749054
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.552879
$ws.Range("N5").Value = 4.658637000000001
$ws.Range("O5").Value = 0.4544999291709252
$ws.Range("P5").Value = 0.4544999291709252
$ws.Range("Q5").Value = 107.2273006979657
$ws.Range("R5").Value = 965.0457062816911
$ws.Range("S5").Value = 0.4051611675824513
$ws.Range("T5").Value = 0.4051611675824512

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Ccl11"
$ws.Range("C6").Value = "Cxcr3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.558577
$ws.Range("H6").Value = 4.675731
$ws.Range("I6").Value = 0.02012122885778811
$ws.Range("J6").Value = 0.02012122885778811
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 1.863797
$ws.Range("N6").Value = 5.591391
$ws.Range("O6").Value = 0.5455000708290748
$ws.Range("P6").Value = 0.5455000708290748
$ws.Range("Q6").Value = 2.904871136869
$ws.Range("R6").Value = 26.143840231821
$ws.Range("S6").Value = 0.01097613176709144
$ws.Range("T6").Value = 0.01097613176709144

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ccl11"
$ws.Range("C7").Value = "Cxcr3"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.558577
$ws.Range("H7").Value = 4.675731
$ws.Range("I7").Value = 0.02012122885778811
$ws.Range("J7").Value = 0.02012122885778811
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.552879
$ws.Range("N7").Value = 4.658637000000001
$ws.Range("O7").Value = 0.4544999291709252
$ws.Range("P7").Value = 0.4544999291709252
$ws.Range("Q7").Value = 2.420281493183
$ws.Range("R7").Value = 21.782533438647
$ws.Range("S7").Value = 0.009145097090696672
$ws.Range("T7").Value = 0.009145097090696672

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ccl11"
$ws.Range("C8").Value = "Cxcr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.007504
$ws.Range("H8").Value = 18.022512
$ws.Range("I8").Value = 0.07755687582203348
$ws.Range("J8").Value = 0.07755687582203348
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 1.863797
$ws.Range("N8").Value = 5.591391
$ws.Range("O8").Value = 0.5455000708290748
$ws.Range("P8").Value = 0.5455000708290748
$ws.Range("Q8").Value = 11.196767932688
$ws.Range("R8").Value = 100.770911394192
$ws.Range("S8").Value = 0.04230728125420102
$ws.Range("T8").Value = 0.04230728125420102

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ccl11"
$ws.Range("C9").Value = "Cxcr3"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.007504
$ws.Range("H9").Value = 18.022512
$ws.Range("I9").Value = 0.07755687582203348
$ws.Range("J9").Value = 0.07755687582203348
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.552879
$ws.Range("N9").Value = 4.658637000000001
$ws.Range("O9").Value = 0.4544999291709252
$ws.Range("P9").Value = 0.4544999291709252
$ws.Range("Q9").Value = 9.328926804016001
$ws.Range("R9").Value = 83.960341236144
$ws.Range("S9").Value = 0.03524959456783246
$ws.Range("T9").Value = 0.03524959456783246

